# OW-248 externalized the graph data files into acuo-data
#
# The placeholder account id "acc1" (shown under "Position Account ID" /
# "FCM Name") is replaced everywhere by the real, externalized account id
# "ACUOSG8745". Both sheets reference the same account, so both B2 cells
# are updated. Row 2 on IRS-Cleared had a tall, explicitly-fixed row
# height to fit the old label; with the new value that custom height is
# no longer needed, so the row reverts to the sheet's default auto
# height. Finally the workbook's active sheet / selection is left on
# IRS-Cleared!D15 (where the edit was reviewed), with IRS-Bilateral's own
# cursor parked at B2.

$wb = $excel.ActiveWorkbook

$wsCleared   = $wb.Worksheets.Item("IRS-Cleared")
$wsBilateral = $wb.Worksheets.Item("IRS-Bilateral")

# Update the account id value on both sheets (they share the same text).
$wsCleared.Range("B2").Value   = "ACUOSG8745"
$wsBilateral.Range("B2").Value = "ACUOSG8745"

# Row 2 on IRS-Cleared no longer needs its tall custom row height now
# that the cell content changed - let it go back to auto height.
$wsCleared.Rows(2).AutoFit()

# Park IRS-Bilateral's cursor at B2 ...
$null = $wsBilateral.Activate()
$null = $wsBilateral.Range("B2").Select()

# ... and make IRS-Cleared the active sheet with D15 selected.
$null = $wsCleared.Activate()
$null = $wsCleared.Range("D15").Select()
